$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header texts:
#    "<field>_old" -> "<field>_FV2404"
#    "<field>_new" -> "<field>_FV2410"
#    (the "diff" header in column K stays untouched)
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2404"
    "B1" = "Segmentgruppe_FV2404"
    "C1" = "Segment_FV2404"
    "D1" = "Datenelement_FV2404"
    "E1" = "Segment ID_FV2404"
    "F1" = "Code_FV2404"
    "G1" = "Qualifier_FV2404"
    "H1" = "Beschreibung_FV2404"
    "I1" = "Bedingungsausdruck_FV2404"
    "J1" = "Bedingung_FV2404"
    "L1" = "Segmentname_FV2410"
    "M1" = "Segmentgruppe_FV2410"
    "N1" = "Segment_FV2410"
    "O1" = "Datenelement_FV2410"
    "P1" = "Segment ID_FV2410"
    "Q1" = "Code_FV2410"
    "R1" = "Qualifier_FV2410"
    "S1" = "Beschreibung_FV2410"
    "T1" = "Bedingungsausdruck_FV2410"
    "U1" = "Bedingung_FV2410"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U76 into an Excel Table ("Table1") with an AutoFilter and
#    banded rows, without pulling in a built-in table style (so that the
#    existing header cell formatting / styles.xml stay untouched).
#
#    Adding a ListObject on top of a range whose header already carries
#    direct formatting makes Excel capture that formatting into a new
#    dxf (headerRowDxfId) - to avoid that we temporarily strip the header
#    formatting, create the table, and then restore the original
#    formatting via copy/paste (which reuses the existing style record
#    instead of creating new font/fill/style entries).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$tempRange = $ws.Range("A200:U200")

$tempRange.Value = "x"
$headerRange.Copy()
$tempRange.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$tempRange.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$tempRange.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.SplitColumn = 0
$win.Split = $false
$win.FreezePanes = $true

Write-Host "Edit complete"
